$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Simple price-only updates (column D)
Set-TextValue "D2"  "246.29"
Set-TextValue "D4"  "5.461"
Set-TextValue "D5"  "0.05628"
Set-TextValue "D6"  "6.468"
Set-TextValue "D7"  "0.8057"
Set-TextValue "D8"  "1.045"
Set-TextValue "D9"  "0.1436"
Set-TextValue "D10" "0.07324"
Set-TextValue "D11" "0.03200"
Set-TextValue "D12" "0.02937"
Set-TextValue "D13" "0.09265"
Set-TextValue "D14" "0.001675"
Set-TextValue "D15" "3.196"
Set-TextValue "D16" "0.04726"

# Row 17: price + Volume(1h) label (drop "Worstin24h")
Set-TextValue "D17" "0.0005827"
$ws.Range("E17").Value = "16OneONE"

# Row 18: price only
Set-TextValue "D18" "0.006347"

# Rows 19-27: coins shifted up one slot with new pricing/labels
$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D19" "0.001056"
$ws.Range("E19").Value = "18BitKanKAN"

$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D20" "0.004103"
$ws.Range("E20").Value = "19HotbitTokenHTB"

$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue "D21" "0.0001504"
$ws.Range("E21").Value = "20NitroExNTX"

$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D22" "3.980"
$ws.Range("E22").Value = "21LEOLEO"

$ws.Range("B23").Value = "GateToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D23" "3.387"
$ws.Range("E23").Value = "22GateTokenGT"

$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D24" "2.085"
$ws.Range("E24").Value = "23BTSETokenBTSE"

$ws.Range("B25").Value = "BitpandaEcosystemToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D25" "0.3267"
$ws.Range("E25").Value = "24BitpandaEcosystemTokenBEST"

$ws.Range("B26").Value = "ProBitToken"
$ws.Range("C26").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D26" "0.1314"
$ws.Range("E26").Value = "25ProBitTokenPROBBestin24h"

$ws.Range("B27").Value = "UpBots"
$ws.Range("C27").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextValue "D27" "0.0002908"
$ws.Range("E27").Value = "26UpBotsUBXT"

# Remaining price-only updates (column D) further down the sheet
Set-TextValue "D40" "0.04158"
Set-TextValue "D41" "0.006910"
Set-TextValue "D43" "0.1037"
Set-TextValue "D44" "0.009019"
Set-TextValue "D45" "0.00005654"
Set-TextValue "D47" "0.6820"

# Row 48: price + Volume(1h) label (add "Worstin24h")
Set-TextValue "D48" "0.01804"
$ws.Range("E48").Value = "47BOLOBOLOWorstin24h"
